# Updated cryptos list (Price + Volume(1h) columns) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "19.762.10"
$ws.Range("E2").Value = "  -8.69%  "
$ws.Range("D3").Value = "1.389.20"
$ws.Range("E3").Value = "  -9.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "268.50"
$ws.Range("E6").Value = "  -7.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3629"
$ws.Range("E7").Value = "  -7.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3030"
$ws.Range("E8").Value = "  -4.38%  "
$ws.Range("E9").Value = "  -7.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9811"
$ws.Range("E10").Value = "  -6.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06405"
$ws.Range("E11").Value = "  -10.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.287"
$ws.Range("E13").Value = "  -6.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.77"
$ws.Range("E14").Value = "  -9.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.063"
$ws.Range("E15").Value = "  -7.95%  "
$ws.Range("D16").Value = "1.389.36"
$ws.Range("E16").Value = "  -9.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009922"
$ws.Range("E17").Value = "  -9.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05620"
$ws.Range("E18").Value = "  -14.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  -16.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.515"
$ws.Range("E21").Value = "  -9.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.38"
$ws.Range("E22").Value = "  -6.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.242"
$ws.Range("E24").Value = "  -4.50%  "
$ws.Range("D25").Value = "19.789.12"
$ws.Range("E25").Value = "  -8.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.182"
$ws.Range("E26").Value = "  -6.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "136.39"
$ws.Range("E27").Value = "  -8.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.62"
$ws.Range("E28").Value = "  -9.03%  "
$ws.Range("D29").Value = "1.545.85"
$ws.Range("E29").Value = "  -9.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "107.77"
$ws.Range("E30").Value = "  -7.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.856"
$ws.Range("E31").Value = "  -20.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.208"
$ws.Range("E32").Value = "  -13.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8020"
$ws.Range("E33").Value = "  -13.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07562"
$ws.Range("E34").Value = "  -6.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.137"
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9993"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05600"
$ws.Range("E37").Value = "  -6.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.679"
$ws.Range("E38").Value = "  -9.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02018"
$ws.Range("E39").Value = "  -8.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1889"
$ws.Range("E40").Value = "  -6.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.12"
$ws.Range("E41").Value = "  -7.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.301"
$ws.Range("E42").Value = "  -10.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.052"
$ws.Range("E43").Value = "  -10.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5195"
$ws.Range("E44").Value = "  -9.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.05"
$ws.Range("E45").Value = "  -7.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.460"
$ws.Range("E46").Value = "  -6.78%  "
$ws.Range("E47").Value = "  -9.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.02"
$ws.Range("E48").Value = "  -5.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.736"
$ws.Range("E49").Value = "  -7.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.032"
$ws.Range("E50").Value = "  -11.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9973"
$ws.Range("E51").Value = "  -0.36%  "
